# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Mapping of row -> new F value (only rows whose value changed)
$updates = @{
    2  = 359
    3  = 353
    4  = 1860
    8  = 741
    10 = 357
    11 = 4442
    13 = 335
    14 = 1228
    15 = 519
    16 = 50
    17 = 796
    19 = 428
    20 = 54
    21 = 209
    22 = 19
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
